$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rectangle 33" - heading above "customer register" box:
# Vaske kunderegister -> Tilby kunde Digipost-faktura
$shpWash = $s.Shapes.Item("Rectangle 33")
$shpWash.TextFrame.TextRange.Text = "Tilby kunde Digipost-faktura"

# "Rectangle 39" - the "1) Initiell vask" label box:
#  - resized/repositioned slightly
#  - text changed to "1) Identifisering" (kept as 3 runs: "1" / ") " / "Identifisering")
$shpInit = $s.Shapes.Item("Rectangle 39")
$shpInit.Left = 101.77307086614174
$shpInit.Top = 117.25165354330709
$shpInit.Width = 110.90748231496063
$shpInit.Height = 38.769213598425196

$tr = $shpInit.TextFrame.TextRange
# Edit from the end of the string backwards so earlier character offsets
# stay valid as the text length changes.
$tr.Characters(12, 5).Text = "Identifisering"
$tr.Characters(4, 8).Text = ") "
$tr.Characters(1, 3).Text = "1"
